$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "data\output\output_subpreg\CE\00200\4000082"
$ws.Range("B3").Value = "No existen archivos disponibles para estudiante serie 4000082"
$ws.Range("C3").Value = "Estudiante"
